$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "69+25=", "45-29=", "90-88=", "33-28=", "23-6=",
    "23+29=", "27+15=", "35+49=", "64-7=", "85-6=",
    "72-39=", "72-19=", "18+27=", "83-15=", "72-49=",
    "60-48=", "58-19=", "97-69=", "27+4=", "71-58=",
    "10-3=", "14+27=", "90-74=", "93-49=", "38+34=",
    "71-68=", "38+6=", "15+17=", "24-19=", "95-38=",
    "34+8=", "76-48=", "14+28=", "17+67=", "81-7=",
    "70-1=", "90-5=", "33+48=", "35-18=", "80-59=",
    "2+19=", "82-63=", "29+44=", "23+68=", "23+29=",
    "6+79=", "63-17=", "39+2=", "93-55=", "83-48=",
    "71-68=", "17+78=", "59+8=", "18+9=", "33-17=",
    "11-3=", "54-15=", "2+49=", "73-8=", "5+16=",
    "27-18=", "90-66=", "67+6=", "68+26=", "47+7=",
    "56-28=", "81-57=", "74-5=", "15+49=", "58+34=",
    "81-25=", "77+16=", "83-59=", "72-27=", "90-83=",
    "48-39=", "60-4=", "31-5=", "33+28=", "84-28=",
    "75+19=", "80-32=", "85-16=", "7+86=", "23+19=",
    "14-7=", "58+26=", "80-56=", "83-79=", "42-36=",
    "32-4=", "40-11=", "28+14=", "26+5=", "62-53=",
    "27+59=", "94-55=", "45-29=", "34+9=", "8+37="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
